$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2: Natalie's - Orange Juice -> Quantity 4->2, Total Cost 95.00->47.50
Set-TextValue $ws.Range("C2") "2"
Set-TextValue $ws.Range("E2") "47.50"

# Row 4: Natalie's - Orange Pineapple -> Quantity 3->2, Total Cost 39.00->26.00
Set-TextValue $ws.Range("C4") "2"
Set-TextValue $ws.Range("E4") "26.00"

# Row 5: Natalie's - Orange Mango -> Quantity 3->2, Total Cost 39.00->26.00
Set-TextValue $ws.Range("C5") "2"
Set-TextValue $ws.Range("E5") "26.00"

# Row 7: Natalie's - Honey Tangerine -> Quantity 2->1, Total Cost 28.00->14.00
Set-TextValue $ws.Range("C7") "1"
Set-TextValue $ws.Range("E7") "14.00"

# New row 9: Natalie's - Lemonade
Set-TextValue $ws.Range("A9") "004011"
Set-TextValue $ws.Range("B9") "Natalie's - Lemonade"
Set-TextValue $ws.Range("C9") "1"
Set-TextValue $ws.Range("D9") "9.25"
Set-TextValue $ws.Range("E9") "9.25"

# New row 10: Natalie's - Strawberry Lemonade
Set-TextValue $ws.Range("A10") "004014"
Set-TextValue $ws.Range("B10") "Natalie's - Strawberry Lemonade"
Set-TextValue $ws.Range("C10") "1"
Set-TextValue $ws.Range("D10") "13.90"
Set-TextValue $ws.Range("E10") "13.90"
